$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new Task rows 19-35 first (appends new shared strings before the row 13 edits)
$ws.Range("A19").Value = "Task 19:"
$ws.Range("A20").Value = "Task 20:"
$ws.Range("A21").Value = "Task 21:"
$ws.Range("A22").Value = "Task 22:"
$ws.Range("A23").Value = "Task 23:"
$ws.Range("A24").Value = "Task 24:"
$ws.Range("A25").Value = "Task 25:"
$ws.Range("A26").Value = "Task 26:"
$ws.Range("A27").Value = "Task 27: "
$ws.Range("A28").Value = "Task 28:"
$ws.Range("A29").Value = "Task 29:"
$ws.Range("A30").Value = "Task 30:"
$ws.Range("A31").Value = "Task 31:"
$ws.Range("A32").Value = "Task 32:"
$ws.Range("A33").Value = "Task 33:"
$ws.Range("A34").Value = "Task 34:"
$ws.Range("A35").Value = "Task 35:"

# Rename Task 13 and record the new password layout bug note
$ws.Range("A13").Value = "Task 13: Modal form sizes for different screens"
$ws.Range("C13").Value = "annoying bug where password interrups the layout of password form"

# Reflect the last selected cell in the sheet view
$ws.Range("C14").Select()
